# Update skill-level values on the Individual Skill Audit sheet
# (commit: "Added Sushant's Belbin's Analysis and updated Team Belbin's Analysis")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Skill Level (1-5) column updates
$ws.Range("C13").Value = 4
$ws.Range("C15").Value = 3
$ws.Range("C22").Value = 3

# Reflect the cursor/scroll position left behind by the editing session
$ws.Range("C23").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
